$d = $word.ActiveDocument

$pairs = @(
    @("721×6=4326", "148×3=444"),
    @("628×5=3140", "715×2=1430"),
    @("592×4=2368", "266×7=1862"),
    @("898×5=4490", "671×7=4697"),
    @("851×5=4255", "599×9=5391"),
    @("721×4=2884", "402×5=2010"),
    @("153×9=1377", "593×2=1186"),
    @("607×4=2428", "166×5=830"),
    @("535×9=4815", "207×7=1449"),
    @("659×4=2636", "343×8=2744"),
    @("738×5=3690", "380×2=760"),
    @("485×2=970", "488×9=4392"),
    @("999×9=8991", "521×6=3126"),
    @("845×3=2535", "509×4=2036"),
    @("499×8=3992", "227×6=1362"),
    @("177×9=1593", "928×2=1856"),
    @("427×8=3416", "785×7=5495"),
    @("977×3=2931", "196×6=1176"),
    @("780×5=3900", "934×7=6538"),
    @("551×2=1102", "678×7=4746"),
    @("323×2=646", "978×2=1956"),
    @("458×6=2748", "950×4=3800"),
    @("446×9=4014", "901×5=4505"),
    @("869×4=3476", "851×4=3404"),
    @("701×3=2103", "996×9=8964")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$d.Save()
